$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 428.77777
$ws.Range("I9").Value = 393.85715
$ws.Range("J9").Value = 551
$ws.Range("K9").Value = 393.85715
$ws.Range("L9").Value = 551
$ws.Range("M9").Value = -224.85715
$ws.Range("N9").Value = -889
$ws.Range("H15").Value = 490
$ws.Range("I15").Value = 490
$ws.Range("K15").Value = 1470
$ws.Range("M15").Value = -1301
$ws.Range("H43").Value = 6918.75
$ws.Range("I43").Value = 6670
$ws.Range("J43").Value = 7333.3335
$ws.Range("K43").Value = 6670
$ws.Range("L43").Value = 7333.3335
$ws.Range("M43").Value = -6601
$ws.Range("N43").Value = -7471.3335
$ws.Range("H70").Value = 3714.2144
$ws.Range("I70").Value = 3500
$ws.Range("J70").Value = 4249.75
$ws.Range("K70").Value = 10500
$ws.Range("L70").Value = 12749.25
$ws.Range("M70").Value = -10230
$ws.Range("N70").Value = -13289.25
$ws.Range("H73").Value = 3714.2144
$ws.Range("I73").Value = 3500
$ws.Range("J73").Value = 4249.75
$ws.Range("K73").Value = 10500
$ws.Range("L73").Value = 12749.25
$ws.Range("M73").Value = -9564
$ws.Range("N73").Value = -14621.25
$ws.Range("H112").Value = 1788.8889
$ws.Range("H132").Value = 2859.2
$ws.Range("J132").Value = 1338
$ws.Range("L132").Value = 4014
$ws.Range("N132").Value = -9074
$ws.Range("H137").Value = 1474.75
$ws.Range("I137").Value = 966.3333
$ws.Range("J137").Value = 3000
$ws.Range("K137").Value = 2898.9999
$ws.Range("L137").Value = 9000
$ws.Range("M137").Value = -348.9998999999998
$ws.Range("N137").Value = -14100

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2801
$ws.Range("I45").Value = 3068
$ws.Range("K45").Value = 3068
$ws.Range("M45").Value = -2691
$ws.Range("H74").Value = 4000
$ws.Range("J74").Value = 4000
$ws.Range("L74").Value = 4000
$ws.Range("N74").Value = -5748
$ws.Range("H77").Value = 4000
$ws.Range("J77").Value = 4000
$ws.Range("L77").Value = 20000
$ws.Range("N77").Value = -28736
$ws.Range("H110").Value = 50000640
$ws.Range("I110").Value = 1100
$ws.Range("K110").Value = 1100
$ws.Range("M110").Value = 945
$ws.Range("H132").Value = 3006
$ws.Range("I132").Value = 3006
$ws.Range("K132").Value = 9018
$ws.Range("M132").Value = -6488
$ws.Range("H139").Value = 99748.75
$ws.Range("J139").Value = 99748.75
$ws.Range("L139").Value = 99748.75
$ws.Range("N139").Value = -110028.75

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4997.5
$ws.Range("J134").Value = 6000
$ws.Range("L134").Value = 18000
$ws.Range("N134").Value = -23070

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 44888.777
$ws.Range("J74").Value = 47599.8
$ws.Range("L74").Value = 47599.8
$ws.Range("N74").Value = -49347.8
$ws.Range("H77").Value = 44888.777
$ws.Range("J77").Value = 47599.8
$ws.Range("L77").Value = 142799.4
$ws.Range("N77").Value = -151535.4
$ws.Range("H141").Value = 372097
$ws.Range("J141").Value = 372097
$ws.Range("L141").Value = 372097
$ws.Range("N141").Value = -382457

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 843.75
$ws.Range("I17").Value = 1000
$ws.Range("K17").Value = 3000
$ws.Range("M17").Value = -2831
$ws.Range("H34").Value = 2533.3333
$ws.Range("I34").Value = 134
$ws.Range("J34").Value = 4932.6665
$ws.Range("K34").Value = 402
$ws.Range("L34").Value = 14797.9995
$ws.Range("M34").Value = -318
$ws.Range("N34").Value = -14965.9995
$ws.Range("H39").Value = 1003
$ws.Range("I39").Value = 1003
$ws.Range("K39").Value = 3009
$ws.Range("M39").Value = -2715
$ws.Range("H55").Value = 3844.7163
$ws.Range("J55").Value = 3960.1904
$ws.Range("L55").Value = 11880.5712
$ws.Range("N55").Value = -12234.5712
$ws.Range("H64").Value = 3866.5
$ws.Range("J64").Value = 4131.4287
$ws.Range("L64").Value = 12394.2861
$ws.Range("N64").Value = -12934.2861
$ws.Range("H67").Value = 3866.5
$ws.Range("J67").Value = 4131.4287
$ws.Range("L67").Value = 12394.2861
$ws.Range("N67").Value = -14266.2861
$ws.Range("H92").Value = 1172
$ws.Range("I92").Value = 726.3333
$ws.Range("J92").Value = 1506.25
$ws.Range("K92").Value = 2178.9999
$ws.Range("L92").Value = 4518.75
$ws.Range("M92").Value = -930.9998999999998
$ws.Range("N92").Value = -7014.75
$ws.Range("H116").Value = 445.25
$ws.Range("I116").Value = 445.25
$ws.Range("K116").Value = 1335.75
$ws.Range("M116").Value = 2106.25
$ws.Range("H129").Value = 2779.25
$ws.Range("J129").Value = 2483.3333
$ws.Range("L129").Value = 7449.999899999999
$ws.Range("N129").Value = -17449.9999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H75").Value = 51153.848
$ws.Range("J75").Value = 51153.848
$ws.Range("L75").Value = 51153.848
$ws.Range("N75").Value = -52901.848
$ws.Range("H78").Value = 51153.848
$ws.Range("J78").Value = 51153.848
$ws.Range("L78").Value = 153461.544
$ws.Range("N78").Value = -162197.544
$ws.Range("H132").Value = 4604.6
$ws.Range("I132").Value = 4604.6
$ws.Range("K132").Value = 13813.8
$ws.Range("M132").Value = -11283.8

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").Value = $null
$ws.Range("H122").Value = 5091.9
$ws.Range("I122").Value = 4718.5
$ws.Range("J122").Value = 8452.5
$ws.Range("K122").Value = 14155.5
$ws.Range("L122").Value = 25357.5
$ws.Range("M122").Value = -11705.5
$ws.Range("N122").Value = -30257.5
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").Value = $null
$ws.Range("H132").Value = 2301
$ws.Range("I132").Value = 1649.6666
$ws.Range("K132").Value = 4948.9998
$ws.Range("M132").Value = -2418.9998

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 96999.5
$ws.Range("J46").Value = 96999.5
$ws.Range("L46").Value = 96999.5
$ws.Range("N46").Value = -97461.5
$ws.Range("H107").Value = 494.33334
$ws.Range("I107").Value = 350
$ws.Range("K107").Value = 1050
$ws.Range("M107").Value = 870
$ws.Range("H125").Value = 40715
$ws.Range("J125").Value = 40715
$ws.Range("L125").Value = 40715
$ws.Range("N125").Value = -50555
$ws.Range("H134").Value = 96999.5
$ws.Range("J134").Value = 96999.5
$ws.Range("L134").Value = 290998.5
$ws.Range("N134").Value = -296068.5
